$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 0.01608704277522737
$ws.Cells.Item(2, 3).Value = 0.01608704277522737
$ws.Cells.Item(2, 4).Value = 99
$ws.Cells.Item(2, 5).Value = 99
$ws.Cells.Item(2, 6).Value = 187
$ws.Cells.Item(2, 7).Value = 187

$ws.Cells.Item(3, 2).Value = 0.4504501307805319
$ws.Cells.Item(3, 3).Value = 0.4504501307805319
$ws.Cells.Item(3, 4).Value = 55
$ws.Cells.Item(3, 5).Value = 55
$ws.Cells.Item(3, 6).Value = 76
$ws.Cells.Item(3, 7).Value = 76

$ws.Cells.Item(4, 2).Value = 0.599326575232672
$ws.Cells.Item(4, 3).Value = 0.599326575232672
$ws.Cells.Item(4, 4).Value = 27
$ws.Cells.Item(4, 5).Value = 27
$ws.Cells.Item(4, 6).Value = 49
$ws.Cells.Item(4, 7).Value = 49

$ws.Cells.Item(5, 2).Value = 0.003435140215643053
$ws.Cells.Item(5, 3).Value = 0.002887058754728575
$ws.Cells.Item(5, 4).Value = 267
$ws.Cells.Item(5, 5).Value = 256
$ws.Cells.Item(5, 6).Value = 271
$ws.Cells.Item(5, 7).Value = 271

$ws.Cells.Item(6, 2).Value = 0.9391719980261437
$ws.Cells.Item(6, 3).Value = 0.9391719980261437
$ws.Cells.Item(6, 4).Value = 7
$ws.Cells.Item(6, 5).Value = 7
$ws.Cells.Item(6, 6).Value = 9
$ws.Cells.Item(6, 7).Value = 9

$ws.Cells.Item(7, 2).Value = 0.00003406597716935116
$ws.Cells.Item(7, 3).Value = 0.00003276167794906654
$ws.Cells.Item(7, 4).Value = 884
$ws.Cells.Item(7, 5).Value = 877
$ws.Cells.Item(7, 6).Value = 898
$ws.Cells.Item(7, 7).Value = 898

$ws.Cells.Item(8, 2).Value = 0.06788442715084926
$ws.Cells.Item(8, 3).Value = 0.06319555246298783
$ws.Cells.Item(8, 4).Value = 169
$ws.Cells.Item(8, 5).Value = 162
$ws.Cells.Item(8, 6).Value = 172
$ws.Cells.Item(8, 7).Value = 172

$ws.Cells.Item(9, 2).Value = 0.03587035260617643
$ws.Cells.Item(9, 3).Value = 0.03034398290622116
$ws.Cells.Item(9, 4).Value = 297
$ws.Cells.Item(9, 5).Value = 267
$ws.Cells.Item(9, 6).Value = 306
$ws.Cells.Item(9, 7).Value = 306

$ws.Cells.Item(10, 2).Value = 0.3269407760033544
$ws.Cells.Item(10, 3).Value = 0.2897441732233381
$ws.Cells.Item(10, 4).Value = 102
$ws.Cells.Item(10, 5).Value = 74
$ws.Cells.Item(10, 6).Value = 125
$ws.Cells.Item(10, 7).Value = 122

$ws.Cells.Item(11, 2).Value = 0.00000000009960571378141858
$ws.Cells.Item(11, 3).Value = 0.00000000008083473722803328
$ws.Cells.Item(11, 4).Value = 1443
$ws.Cells.Item(11, 5).Value = 1422
$ws.Cells.Item(11, 6).Value = 1488
$ws.Cells.Item(11, 7).Value = 1477

$ws.Cells.Item(12, 2).Value = 0.010418072253105
$ws.Cells.Item(12, 3).Value = 0.010418072253105
$ws.Cells.Item(12, 4).Value = 217
$ws.Cells.Item(12, 5).Value = 217
$ws.Cells.Item(12, 6).Value = 312
$ws.Cells.Item(12, 7).Value = 312

$ws.Cells.Item(13, 2).Value = 0.5035560325870316
$ws.Cells.Item(13, 3).Value = 0.5004452583181059
$ws.Cells.Item(13, 4).Value = 90
$ws.Cells.Item(13, 5).Value = 86
$ws.Cells.Item(13, 6).Value = 92
$ws.Cells.Item(13, 7).Value = 92

$ws.Cells.Item(14, 2).Value = 0.7688953153587761
$ws.Cells.Item(14, 3).Value = 0.7617814239549737
$ws.Cells.Item(14, 4).Value = 30
$ws.Cells.Item(14, 5).Value = 25
$ws.Cells.Item(14, 6).Value = 36
$ws.Cells.Item(14, 7).Value = 36

$ws.Cells.Item(15, 2).Value = 0.7676979607076297
$ws.Cells.Item(15, 3).Value = 0.7652295892236318
$ws.Cells.Item(15, 4).Value = 28
$ws.Cells.Item(15, 5).Value = 27
$ws.Cells.Item(15, 6).Value = 32
$ws.Cells.Item(15, 7).Value = 31

$ws.Cells.Item(16, 2).Value = 0.6266129947137067
$ws.Cells.Item(16, 3).Value = 0.6266129947137067
$ws.Cells.Item(16, 4).Value = 62
$ws.Cells.Item(16, 5).Value = 62
$ws.Cells.Item(16, 6).Value = 65
$ws.Cells.Item(16, 7).Value = 65

$ws.Cells.Item(17, 2).Value = 0.4686914267282534
$ws.Cells.Item(17, 3).Value = 0.4686914267282534
$ws.Cells.Item(17, 4).Value = 36
$ws.Cells.Item(17, 5).Value = 36
$ws.Cells.Item(17, 6).Value = 63
$ws.Cells.Item(17, 7).Value = 63

